# Add a new "target_type" enum row: enemy_all / TARGET_TYPE_ENEMY_ALL / 4 / 적 전체
# This is inserted right after the existing TARGET_TYPE_ALLY_HP_LOWEST row (row 70),
# which pushes every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71 (shifts old row 71.. down to 72..)
$ws.Rows("71:71").Insert()

# Fill in the new row's data
$ws.Range("A71").Value = "target_type"
$ws.Range("B71").Value = "enemy_all"
$ws.Range("C71").Formula = "=UPPER(A71)&""_""&UPPER(B71)"
$ws.Range("D71").Value = 4
$ws.Range("E71").Value = "적 전체"

# Restore/update the view: the frozen pane now scrolls further down and the
# previously selected cell (which has shifted from row 60 to row 72) is
# re-selected.
$app = $ws.Application
$win = $app.ActiveWindow
$win.ScrollRow = 48
[void]$ws.Range("E72").Select()

Write-Output "done"
